$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1) - F column "想去人数" updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 93
$ws1.Cells.Item(3, 6).Value = 170
$ws1.Cells.Item(4, 6).Value = 405
$ws1.Cells.Item(5, 6).Value = 183
$ws1.Cells.Item(7, 6).Value = 1068
$ws1.Cells.Item(8, 6).Value = 354
$ws1.Cells.Item(9, 6).Value = 180
$ws1.Cells.Item(10, 6).Value = 48
$ws1.Cells.Item(12, 6).Value = 365
$ws1.Cells.Item(13, 6).Value = 360
$ws1.Cells.Item(14, 6).Value = 772
$ws1.Cells.Item(15, 6).Value = 142
$ws1.Cells.Item(16, 6).Value = 707
$ws1.Cells.Item(18, 6).Value = 71
$ws1.Cells.Item(19, 6).Value = 984
$ws1.Cells.Item(20, 6).Value = 441
$ws1.Cells.Item(21, 6).Value = 253
$ws1.Cells.Item(23, 6).Value = 370
$ws1.Cells.Item(25, 6).Value = 37

# Sheet 2: 演出 (index 2) - F column updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value = 360
$ws2.Cells.Item(9, 6).Value = 8

# Sheet 4: 全部类型 (index 4) - F column updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(4, 6).Value = 93
$ws4.Cells.Item(5, 6).Value = 170
$ws4.Cells.Item(6, 6).Value = 405
$ws4.Cells.Item(7, 6).Value = 183
$ws4.Cells.Item(9, 6).Value = 1068
$ws4.Cells.Item(10, 6).Value = 355
$ws4.Cells.Item(11, 6).Value = 180
$ws4.Cells.Item(13, 6).Value = 48
$ws4.Cells.Item(14, 6).Value = 360
$ws4.Cells.Item(17, 6).Value = 365
$ws4.Cells.Item(20, 6).Value = 360
$ws4.Cells.Item(21, 6).Value = 772
$ws4.Cells.Item(22, 6).Value = 142
$ws4.Cells.Item(23, 6).Value = 707
$ws4.Cells.Item(25, 6).Value = 71
$ws4.Cells.Item(26, 6).Value = 984
$ws4.Cells.Item(27, 6).Value = 441
$ws4.Cells.Item(29, 6).Value = 8
$ws4.Cells.Item(30, 6).Value = 253
$ws4.Cells.Item(32, 6).Value = 370
$ws4.Cells.Item(36, 6).Value = 37
